# Update "想去人数" (attendee count) figures in column F, as captured by a
# fresh data scrape (gh-pages output regeneration). The same underlying
# event table is duplicated across the "展览" (sheet 1) and "全部类型"
# (sheet 4) worksheets, so every row update is applied to both.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 1832
    "F5"  = 39
    "F7"  = 1503
    "F9"  = 592
    "F11" = 97
    "F12" = 17
    "F16" = 132
    "F19" = 3549
    "F20" = 424
    "F21" = 316
    "F22" = 484
    "F23" = 113
    "F26" = 1306
}

$sheetIndexes = @(1, 4)

foreach ($sheetIndex in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
